$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 60.90318633333334
$ws.Range("H2").Value = 182.709559
$ws.Range("I2").Value = 0.4799022665420342
$ws.Range("J2").Value = 0.4799022665420342
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 3.652835333333333
$ws.Range("N2").Value = 10.958506
$ws.Range("O2").Value = 0.1451640905049821
$ws.Range("P2").Value = 0.1451640905049821
$ws.Range("Q2").Value = 222.4693109509838
$ws.Range("R2").Value = 2002.223798558854
$ws.Range("S2").Value = 0.06966457605385391
$ws.Range("T2").Value = 0.06966457605385391

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 60.90318633333334
$ws.Range("H3").Value = 182.709559
$ws.Range("I3").Value = 0.4799022665420342
$ws.Range("J3").Value = 0.4799022665420342
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 2.483777
$ws.Range("N3").Value = 7.451331
$ws.Range("O3").Value = 0.09870557972652284
$ws.Range("P3").Value = 0.09870557972652286
$ws.Range("Q3").Value = 151.2699334414477
$ws.Range("R3").Value = 1361.429400973029
$ws.Range("S3").Value = 0.04736903143110377
$ws.Range("T3").Value = 0.04736903143110378

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 60.90318633333334
$ws.Range("H4").Value = 182.709559
$ws.Range("I4").Value = 0.4799022665420342
$ws.Range("J4").Value = 0.4799022665420342
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 19.026879
$ws.Range("N4").Value = 57.080637
$ws.Range("O4").Value = 0.7561303297684949
$ws.Range("P4").Value = 0.756130329768495
$ws.Range("Q4").Value = 1158.797557078787
$ws.Range("R4").Value = 10429.17801370908
$ws.Range("S4").Value = 0.3628686590570764
$ws.Range("T4").Value = 0.3628686590570765

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 51.42568199999999
$ws.Range("H5").Value = 154.277046
$ws.Range("I5").Value = 0.405221842009972
$ws.Range("J5").Value = 0.405221842009972
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 3.652835333333333
$ws.Range("N5").Value = 10.958506
$ws.Range("O5").Value = 0.1451640905049821
$ws.Range("P5").Value = 0.1451640905049821
$ws.Range("Q5").Value = 187.849548250364
$ws.Range("R5").Value = 1690.645934253276
$ws.Range("S5").Value = 0.05882366014813114
$ws.Range("T5").Value = 0.05882366014813114

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 51.42568199999999
$ws.Range("H6").Value = 154.277046
$ws.Range("I6").Value = 0.405221842009972
$ws.Range("J6").Value = 0.405221842009972
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 2.483777
$ws.Range("N6").Value = 7.451331
$ws.Range("O6").Value = 0.09870557972652284
$ws.Range("P6").Value = 0.09870557972652286
$ws.Range("Q6").Value = 127.729926160914
$ws.Range("R6").Value = 1149.569335448226
$ws.Range("S6").Value = 0.03999765683344373
$ws.Range("T6").Value = 0.03999765683344374

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 51.42568199999999
$ws.Range("H7").Value = 154.277046
$ws.Range("I7").Value = 0.405221842009972
$ws.Range("J7").Value = 0.405221842009972
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 19.026879
$ws.Range("N7").Value = 57.080637
$ws.Range("O7").Value = 0.7561303297684949
$ws.Range("P7").Value = 0.756130329768495
$ws.Range("Q7").Value = 978.4702289064778
$ws.Range("R7").Value = 8806.2320601583
$ws.Range("S7").Value = 0.306400525028397
$ws.Range("T7").Value = 0.3064005250283971

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 14.57860966666667
$ws.Range("H8").Value = 43.735829
$ws.Range("I8").Value = 0.1148758914479938
$ws.Range("J8").Value = 0.1148758914479938
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 3.652835333333333
$ws.Range("N8").Value = 10.958506
$ws.Range("O8").Value = 0.1451640905049821
$ws.Range("P8").Value = 0.1451640905049821
$ws.Range("Q8").Value = 53.25326050127489
$ws.Range("R8").Value = 479.2793445114739
$ws.Range("S8").Value = 0.01667585430299708
$ws.Range("T8").Value = 0.01667585430299708

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 14.57860966666667
$ws.Range("H9").Value = 43.735829
$ws.Range("I9").Value = 0.1148758914479938
$ws.Range("J9").Value = 0.1148758914479938
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 2.483777
$ws.Range("N9").Value = 7.451331
$ws.Range("O9").Value = 0.09870557972652284
$ws.Range("P9").Value = 0.09870557972652286
$ws.Range("Q9").Value = 36.21001538204433
$ws.Range("R9").Value = 325.8901384383989
$ws.Range("S9").Value = 0.01133889146197534
$ws.Range("T9").Value = 0.01133889146197534

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 14.57860966666667
$ws.Range("H10").Value = 43.735829
$ws.Range("I10").Value = 0.1148758914479938
$ws.Range("J10").Value = 0.1148758914479938
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 19.026879
$ws.Range("N10").Value = 57.080637
$ws.Range("O10").Value = 0.7561303297684949
$ws.Range("P10").Value = 0.756130329768495
$ws.Range("Q10").Value = 277.3854421158969
$ws.Range("R10").Value = 2496.468979043073
$ws.Range("S10").Value = 0.0868611456830214
$ws.Range("T10").Value = 0.08686114568302142

Write-Output "Updated rows 2-10 per Natmi Dr Hou advice"